# Commit: Added custom prompt type: adate
# This adds a new "TEST" form entry (mirroring the existing MIF / MIFVISIT /
# MIFSCAR / MIFVAC pattern) to the "survey" and "choices" sheets, and moves
# the active/selected sheet from "settings" to "framework_translations".

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("settings")
$wsChoices  = $wb.Worksheets.Item("choices")
$wsSurvey   = $wb.Worksheets.Item("survey")
$wsFwkTr    = $wb.Worksheets.Item("framework_translations")

# --- survey sheet: append the "TEST" form branch (rows 20-22) -------------
# Note leading '' is an escaped literal apostrophe so the stored text keeps
# the leading quote character (matches the existing MIF/.. rows) instead of
# being interpreted as an Excel quote-prefix.
$wsSurvey.Cells.Item(20, 1).Value = "TEST"

$wsSurvey.Cells.Item(21, 2).Value = "''?' + odkSurvey.getHashString('TEST')"
$wsSurvey.Cells.Item(21, 5).Value = "external_link"
$wsSurvey.Cells.Item(21, 7).Value = "Open test"

$wsSurvey.Cells.Item(22, 3).Value = "exit section"

# --- choices sheet: append the "TEST" choice row (row 6) ------------------
$wsChoices.Cells.Item(6, 1).Value = "forms"
$wsChoices.Cells.Item(6, 2).Value = "TEST"
$wsChoices.Cells.Item(6, 3).Value = "TEST Form"
$wsChoices.Cells.Item(6, 4).Value = "Test formulario"

# --- selections on the touched sheets --------------------------------------
$wsSettings.Activate() | Out-Null
$wsSettings.Range("A10").Select() | Out-Null

$wsChoices.Activate() | Out-Null
$wsChoices.Range("D7").Select() | Out-Null

$wsSurvey.Activate() | Out-Null
$wsSurvey.Range("A21").Select() | Out-Null

# --- move the active/selected tab from "settings" to "framework_translations"
$wsFwkTr.Activate() | Out-Null
